$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell F1, matching the style of the existing header row (E1)
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

# New data column values (time_taken) for rows 2-8
$times = @(
    "2021-10-05 10:51:10.947971",
    "2021-10-05 10:51:10.947984",
    "2021-10-05 10:51:10.947988",
    "2021-10-05 10:51:10.947991",
    "2021-10-05 10:51:10.947995",
    "2021-10-05 10:51:10.947999",
    "2021-10-05 10:51:10.948002"
)

for ($i = 0; $i -lt $times.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $times[$i]
}
